$wb = $excel.ActiveWorkbook

function Rotate-Rows {
    param(
        [object]$ws,
        [string]$lastCol
    )

    $addr3 = "A3:" + $lastCol + "3"
    $addr4 = "A4:" + $lastCol + "4"
    $addr5 = "A5:" + $lastCol + "5"

    $row3 = $ws.Range($addr3).Value2
    $row4 = $ws.Range($addr4).Value2
    $row5 = $ws.Range($addr5).Value2

    $ws.Range($addr3).Value2 = $row5
    $ws.Range($addr4).Value2 = $row3
    $ws.Range($addr5).Value2 = $row4

    foreach ($hl in $ws.Hyperlinks) {
        $r = $hl.Range.Row
        $c = $hl.Range.Column
        if ($r -ge 3 -and $r -le 5) {
            $hl.TextToDisplay = $ws.Cells.Item($r, $c).Value2
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
Rotate-Rows $wsOverview "C"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Rotate-Rows $wsZhCn "D"

$wsDeDe = $wb.Worksheets.Item("de-de")
Rotate-Rows $wsDeDe "D"

Write-Output "done"
